# Commit: "Download and unarchive scanned pdfs."
#
# Adds four new configuration rows to the "Settings" sheet (Path_DownloadFolder,
# Path_InputFolder, Path_SourceFolder, Path_VotesFolder) describing the folders
# used to download and unarchive the scanned-votes PDF archive, and makes the
# Settings sheet the active/selected sheet (previously Constants was active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert three new blank rows at row 10 - together with the (already blank)
# original row 10 this gives four rows (10-13) of space for the new entries,
# and shifts every following row (CREDENTIALS section, EMAILS section, ...)
# down by three.
$ws.Rows("10:12").Insert()

# Keep the inserted rows' height consistent with the rest of the sheet.
$ws.Rows("10:12").RowHeight = 14.25

# Row 10 - Path_DownloadFolder
$ws.Range("A10").Value = "Path_DownloadFolder"
$ws.Range("B10").Value = "C:\Users\giorgia.cocis\Desktop\SDD and Estimations"
$ws.Range("C10").Value = "The path of the folder where the archive is being downloaded."

# Row 11 - Path_InputFolder
$ws.Range("A11").Value = "Path_InputFolder"
$ws.Range("B11").Value = "C:\Users\giorgia.cocis\Desktop\SDD and Estimations\VoteCounter\Data\Input"
$ws.Range("C11").Value = "The path of the folder where the files are being unarchived."

# Row 12 - Path_SourceFolder
$ws.Range("A12").Value = "Path_SourceFolder"
$ws.Range("B12").Value = "C:\Users\giorgia.cocis\Desktop\SDD and Estimations\ScannedVotes.zip"
$ws.Range("C12").Value = "The path of the folder containing the archive."

# Row 13 - Path_VotesFolder
$ws.Range("A13").Value = "Path_VotesFolder"
$ws.Range("B13").Value = "C:\Users\giorgia.cocis\Desktop\SDD and Estimations\VoteCounter\Data\Input\ScannedVotes"
$ws.Range("C13").Value = "The path of the folder containing the scanned pdfs."

# The author left the workbook with the Settings sheet active (instead of
# Constants) and cell C17 selected there.
$ws.Activate() | Out-Null
$ws.Range("C17").Select() | Out-Null
